$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Row 17 (RSAFS, M/M % Delta) - new "Present" value in F17, values shift
# right (old F/G become G/H... actually present->F, values shift across
# F..J), oldest (former J) drops off.
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 45992
$ws.Range("F17").Value = -0.000160525653495891
$ws.Range("G17").Value = 0.005518082869731433
$ws.Range("H17").Value = -0.001558334425942887
$ws.Range("I17").Value = 0.0006724067240673204
$ws.Range("J17").Value = 0.00545946488174831

# Row 18 (RSAFS, Y/Y % Delta)
$ws.Range("C18").Value = 45992
$ws.Range("F18").Value = 0.02427715536403887
$ws.Range("G18").Value = 0.03263595152369941
$ws.Range("H18").Value = 0.03213514238518121
$ws.Range("I18").Value = 0.04144341481107452
$ws.Range("J18").Value = 0.04972605550048132

# Row 29 (T5YIFR)
$ws.Range("N29").Value = 46062
$ws.Range("Q29").Value = 2.2
$ws.Range("R29").Value = 2.18
$ws.Range("S29").Value = 2.16
$ws.Range("T29").Value = 2.19
$ws.Range("U29").Value = 2.19

# Row 30 (T10YIE)
$ws.Range("N30").Value = 46062
$ws.Range("Q30").Value = 2.35
$ws.Range("R30").Value = 2.34
$ws.Range("S30").Value = 2.32
$ws.Range("T30").Value = 2.35
$ws.Range("U30").Value = 2.36

# Row 31 (ECIWAG, Q/Q % Delta)
$ws.Range("N31").Value = 45931
$ws.Range("Q31").Value = 0.007310493043885868
$ws.Range("R31").Value = 0.007996957929548465
$ws.Range("S31").Value = 0.01027939464493599
$ws.Range("T31").Value = 0.007624633431085215
$ws.Range("U31").Value = 0.009473060982829962

# Row 32 (ECIWAG, Y/Y % Delta)
$ws.Range("N32").Value = 45931
$ws.Range("Q32").Value = 0.03362463343108507
$ws.Range("R32").Value = 0.03584369449378332
$ws.Range("S32").Value = 0.03559665871121723
$ws.Range("T32").Value = 0.03369434416365838
$ws.Range("U32").Value = 0.03710462287104619

# Row 41 (IQ, M/M % Delta)
$ws.Range("N41").Value = 45992
$ws.Range("Q41").Value = 0.003255208333333259
$ws.Range("R41").ClearContents()
$ws.Range("S41").ClearContents()
$ws.Range("T41").Value = 0
$ws.Range("U41").Value = 0.00130890052356003

# Row 42 (IQ, Y/Y % Delta)
$ws.Range("N42").Value = 45992
$ws.Range("Q42").Value = 0.03076923076923073
$ws.Range("R42").Value = 0.03225806451612891
$ws.Range("S42").ClearContents()
$ws.Range("T42").Value = 0.03869653767820766
$ws.Range("U42").Value = 0.03238866396761141

# Row 43 (IR, M/M % Delta)
$ws.Range("N43").Value = 45992
$ws.Range("Q43").Value = 0.00141643059490093
$ws.Range("R43").ClearContents()
$ws.Range("S43").ClearContents()
$ws.Range("T43").Value = -0.001418439716311948
$ws.Range("U43").Value = -0.001416430594900819

# Row 44 (IR, Y/Y % Delta)
$ws.Range("N44").Value = 45992
$ws.Range("Q44").Value = 0
$ws.Range("R44").Value = -0.0007077140835104227
$ws.Range("S44").ClearContents()
$ws.Range("T44").Value = -0.0007097232079488595
$ws.Range("U44").Value = -0.002828854314002869

# Row 47 (DFF)
$ws.Range("N47").Value = 46059

# Row 48 (DGS2)
$ws.Range("N48").Value = 46059
$ws.Range("Q48").Value = 3.5
$ws.Range("R48").Value = 3.47
$ws.Range("U48").Value = 3.57

# Row 49 (DGS5)
$ws.Range("N49").Value = 46059
$ws.Range("Q49").Value = 3.76
$ws.Range("R49").Value = 3.74
$ws.Range("U49").Value = 3.83

# Row 50 (DGS10)
$ws.Range("N50").Value = 46059
$ws.Range("Q50").Value = 4.22
$ws.Range("R50").Value = 4.21
$ws.Range("S50").Value = 4.29
$ws.Range("T50").Value = 4.28
$ws.Range("U50").Value = 4.29

# Row 52 (DBAA)
$ws.Range("N52").Value = 46059
$ws.Range("Q52").Value = 5.87
$ws.Range("R52").Value = 5.88
$ws.Range("S52").Value = 5.93
$ws.Range("T52").Value = 5.91
$ws.Range("U52").Value = 5.9

# ---------------------------------------------------------------------
# Style updates: the "Latest Date" cells for these rows get the
# highlighted (yellow-fill) date style used elsewhere in the sheet
# (e.g. N29). Use copy/paste-special (formats only) so the existing
# shared cell style is reused instead of a new one being created.
# ---------------------------------------------------------------------
$ws.Range("N29").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("N31").PasteSpecial(-4122)
$ws.Range("N32").PasteSpecial(-4122)
$ws.Range("N41").PasteSpecial(-4122)
$ws.Range("N42").PasteSpecial(-4122)
$ws.Range("N43").PasteSpecial(-4122)
$ws.Range("N44").PasteSpecial(-4122)
$excel.CutCopyMode = $false
